# Updates cryptos list cell values (Coin/Link/Price/Volume(1h)) per scraped data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as text (avoids Excel auto-converting numeric-looking
# strings like "1.00" or "198.24" into real numbers), while keeping the cell
# style untouched (matches original unstyled data cells).
function Set-TextCell($Sheet, $Row, $Col, $Text) {
    $cell = $Sheet.Cells.Item($Row, $Col)
    $cell.Value = "'" + $Text
    $cell.Style = "Normal"
}

Set-TextCell $ws 2 4 "67.381.63"
Set-TextCell $ws 2 5 "  -1.77%  "
Set-TextCell $ws 3 4 "3.505.56"
Set-TextCell $ws 3 5 "  -3.91%  "
Set-TextCell $ws 4 4 "1.00"
Set-TextCell $ws 4 5 "  +0.29%  "
Set-TextCell $ws 5 4 "198.24"
Set-TextCell $ws 5 5 "  +0.95%  "
Set-TextCell $ws 6 4 "549.33"
Set-TextCell $ws 6 5 "  -6.02%  "
Set-TextCell $ws 7 4 "3.489.93"
Set-TextCell $ws 7 5 "  -4.03%  "
Set-TextCell $ws 8 4 "0.603"
Set-TextCell $ws 8 5 "  -3.29%  "
Set-TextCell $ws 9 4 "1.00"
Set-TextCell $ws 9 5 "  +0.09%  "
Set-TextCell $ws 10 5 "  -4.54%  "
Set-TextCell $ws 11 4 "62.81"
Set-TextCell $ws 11 5 "  +12.42%  "
Set-TextCell $ws 12 5 "  -8.56%  "
Set-TextCell $ws 13 4 "0.0000269"
Set-TextCell $ws 13 5 "  -10.33%  "
Set-TextCell $ws 14 4 "9.78"
Set-TextCell $ws 14 5 "  -3.92%  "
Set-TextCell $ws 15 4 "4.067.58"
Set-TextCell $ws 15 5 "  -3.58%  "
Set-TextCell $ws 16 4 "3.501.32"
Set-TextCell $ws 16 5 "  -3.85%  "
Set-TextCell $ws 17 5 "  -2.05%  "
Set-TextCell $ws 18 2 "WrappedBTC"
Set-TextCell $ws 18 3 "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextCell $ws 18 4 "67.062.23"
Set-TextCell $ws 18 5 "  -1.89%  "
Set-TextCell $ws 19 2 "Chainlink"
Set-TextCell $ws 19 3 "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextCell $ws 19 4 "18.38"
Set-TextCell $ws 19 5 "  -1.46%  "
Set-TextCell $ws 21 4 "1.02"
Set-TextCell $ws 21 5 "  -5.34%  "
Set-TextCell $ws 22 4 "389.37"
Set-TextCell $ws 22 5 "  -4.12%  "
Set-TextCell $ws 23 4 "3.99"
Set-TextCell $ws 23 5 "  -6.35%  "
Set-TextCell $ws 24 4 "11.76"
Set-TextCell $ws 24 5 "  -7.73%  "
Set-TextCell $ws 25 4 "82.20"
Set-TextCell $ws 25 5 "  -5.02%  "
Set-TextCell $ws 26 4 "12.25"
Set-TextCell $ws 26 5 "  -3.66%  "
Set-TextCell $ws 27 4 "2.78"
Set-TextCell $ws 27 5 "  -6.17%  "
Set-TextCell $ws 28 4 "3.77"
Set-TextCell $ws 28 5 "  -4.15%  "
Set-TextCell $ws 29 4 "8.71"
Set-TextCell $ws 29 5 "  -5.90%  "
Set-TextCell $ws 30 4 "30.99"
Set-TextCell $ws 30 5 "  -3.04%  "
Set-TextCell $ws 31 4 "671.86"
Set-TextCell $ws 31 5 "  -4.07%  "
Set-TextCell $ws 32 5 "  -14.21%  "
Set-TextCell $ws 33 4 "11.69"
Set-TextCell $ws 33 5 "  -4.99%  "
Set-TextCell $ws 34 4 "63.13"
Set-TextCell $ws 34 5 "  -2.96%  "
Set-TextCell $ws 35 5 "  -7.69%  "
Set-TextCell $ws 36 4 "38.42"
Set-TextCell $ws 36 5 "  -10.58%  "
Set-TextCell $ws 37 5 "  +0.18%  "
Set-TextCell $ws 38 4 "0.398"
Set-TextCell $ws 38 5 "  -5.42%  "
Set-TextCell $ws 39 2 "FirstDigitalUSD"
Set-TextCell $ws 39 3 "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextCell $ws 39 4 "0.999"
Set-TextCell $ws 39 5 "  +0.24%  "
Set-TextCell $ws 40 2 "Kaspa"
Set-TextCell $ws 40 3 "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell $ws 40 4 "0.131"
Set-TextCell $ws 40 5 "  -4.00%  "
Set-TextCell $ws 41 4 "3.069.83"
Set-TextCell $ws 41 5 "  -2.75%  "
Set-TextCell $ws 42 5 "  -5.01%  "
Set-TextCell $ws 43 4 "0.0₃0673"
Set-TextCell $ws 43 5 "  -17.61%  "
Set-TextCell $ws 44 4 "2.76"
Set-TextCell $ws 44 5 "  +5.17%  "
Set-TextCell $ws 45 5 "  -13.74%  "
Set-TextCell $ws 46 4 "2.66"
Set-TextCell $ws 46 5 "  -4.26%  "
Set-TextCell $ws 47 5 "  -7.46%  "
Set-TextCell $ws 48 4 "0.126"
Set-TextCell $ws 48 5 "  -4.87%  "
Set-TextCell $ws 49 4 "137.11"
Set-TextCell $ws 49 5 "  -4.15%  "
Set-TextCell $ws 50 2 "THORChain"
Set-TextCell $ws 50 3 "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextCell $ws 50 4 "8.17"
Set-TextCell $ws 50 5 "  -8.72%  "
Set-TextCell $ws 51 2 "ApeXProtocol"
Set-TextCell $ws 51 3 "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextCell $ws 51 4 "2.87"
Set-TextCell $ws 51 5 "  -8.32%  "
